try {
  "hello" | Out-File -FilePath "C:\temp\test.txt"
  Write-Output "wrote C:\temp"
} catch {
  Write-Output ("ERR1: " + $_.Exception.Message)
}
try {
  Get-ChildItem "C:\" | Out-String | Write-Output
} catch {
  Write-Output ("ERR2: " + $_.Exception.Message)
}
try {
  Get-ChildItem "C:\evals\00000000-0000-0000-0000-000000000000" | Out-String | Write-Output
} catch {
  Write-Output ("ERR3: " + $_.Exception.Message)
}
